$d = $word.ActiveDocument

# Change 1: merge the split runs describing office city/state into the final text
$d.Content.Find.Execute(
    "_proc_office_neighborhood_, _proc_office_city_ - _proc_office_state_, e com endereço eletrônico em _proc_office_site_.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "_proc_office_neighborhood_, _proc_office_state_, e com endereço eletrônico em _proc_office_site_.",
    2
)

# Change 3: replace the "_pro_city_, _proc_state_, _proc_date_" run pair with "_proc_today_"
$d.Content.Find.Execute(
    "_pro_city_, _proc_state_, _proc_date_",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "_proc_today_",
    2
)

# Change 2: nudge the rotated stamp image up slightly (vertical offset -88900 EMU -> -86360 EMU)
$shape = $d.Shapes.Item(1)
$shape.Top = -6.8
